$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.582
$ws.Range("E4").Value = 13.418

$ws.Range("E5").Value = 13.41

$ws.Range("A6").Value = -21.108

$ws.Range("A7").Value = -21.038

$ws.Range("E8").Value = 13.718

$ws.Range("A16").Value = -20.83300000000001
$ws.Range("E16").Value = 13.142

$ws.Range("A20").Value = -22.2

$ws.Range("E22").Value = 13.241
